# Update "Översikt LEKEBERG" workbook:
#  1) Every "Förändrad" date in column C (rows 2..211) moves from 2023-09-21
#     (45190) to 2023-09-23 (45192).
#  2) The case "A 50292-2021" (previously row 5) is re-surveyed and now sorts
#     above "A 26878-2020" (previously row 4): the two rows swap places, and
#     the "A 50292-2021" row picks up new counts/species and a new
#     "Knärotsbuffertlänk" (column U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bulk-update the "Förändrad" column for every data row -------------
$ws.Range("C2:C211").Value = 45192

# --- 2) Row 4 becomes the updated "A 50292-2021" case ----------------------
$ws.Range("A4").Value = "A 50292-2021"
$ws.Range("B4").Value = 44459
$ws.Range("C4").Value = 45192
$ws.Range("D4").Value = "ÖREBRO LÄN"
$ws.Range("E4").Value = "LEKEBERG"
$ws.Range("G4").Value = 6.3
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 4
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 8
$ws.Range("R4").Value = "Knärot`r`nRödstrupig piplärka`r`nSpillkråka`r`nTalltita`r`nGrönpyrola`r`nStubbspretmossa`r`nVästlig hakmossa`r`nZontaggsvamp"
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/artfynd/A 50292-2021.xlsx", "A 50292-2021")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/kartor/A 50292-2021.png", "A 50292-2021")'
$ws.Range("U4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/knärot/A 50292-2021.png", "A 50292-2021")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/klagomål/A 50292-2021.docx", "A 50292-2021")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/klagomålsmail/A 50292-2021.docx", "A 50292-2021")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/tillsyn/A 50292-2021.docx", "A 50292-2021")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/tillsynsmail/A 50292-2021.docx", "A 50292-2021")'

# --- 3) Row 5 becomes the (otherwise unchanged) "A 26878-2020" case --------
$ws.Range("A5").Value = "A 26878-2020"
$ws.Range("B5").Value = 43987
$ws.Range("C5").Value = 45192
$ws.Range("D5").Value = "ÖREBRO LÄN"
$ws.Range("E5").Value = "LEKEBERG"
$ws.Range("G5").Value = 7.3
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = "Gräddporing`r`nSpillkråka`r`nUllticka`r`nSprödporing`r`nDropptaggsvamp`r`nVedticka`r`nVågbandad barkbock"
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/artfynd/A 26878-2020.xlsx", "A 26878-2020")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/kartor/A 26878-2020.png", "A 26878-2020")'
$ws.Range("U5").ClearContents()
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/klagomål/A 26878-2020.docx", "A 26878-2020")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/klagomålsmail/A 26878-2020.docx", "A 26878-2020")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/tillsyn/A 26878-2020.docx", "A 26878-2020")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKEBERG/tillsynsmail/A 26878-2020.docx", "A 26878-2020")'

# The sheet keeps every data row at a fixed 15pt height (no autofit); restore
# that after rewriting the wrapped "Artnamn" text in rows 4 and 5, since
# assigning multi-line values otherwise triggers row autofit.
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
